$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Orders")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Orders sheet ---
# F3: "1" -> "120" (must stay stored as text, like the rest of the column)
$cellF3 = $ws1.Range("F3")
$cellF3.NumberFormat = "@"
$cellF3.Value = "120"
$cellF3.Style = "Normal"

# New row 4 with a single text value in C4
$cellC4 = $ws1.Range("C4")
$cellC4.NumberFormat = "@"
$cellC4.Value = "14_波浪浅紫洋桔梗_Wavy Light Purple Lisianthus_Eustoma grandiflorum (Raf.) Shinners"
$cellC4.Style = "Normal"

# --- Summary sheet ---
# G2: "01801" -> "01801200" (must stay text to preserve the leading zero)
$cellG2 = $ws2.Range("G2")
$cellG2.NumberFormat = "@"
$cellG2.Value = "01801200"
$cellG2.Style = "Normal"
